# Financial data refresh for NJDCY yearly financials (Income Statement,
# Balance Sheet, Cash Flow Statement) - "Doing Updates for Financials"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Total Revenue
$ws.Range("D8").Value = 13452300
$ws.Range("E8").Value = 10841800
$ws.Range("F8").Value = 10651700
$ws.Range("G8").Value = 9296600
$ws.Range("H8").Value = 7911000
$ws.Range("I8").Value = 6411800
$ws.Range("J8").Value = 6168200

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 10223600
$ws.Range("E9").Value = 8259700
$ws.Range("F9").Value = 8226000
$ws.Range("G9").Value = 7109800
$ws.Range("H9").Value = 6101100
$ws.Range("I9").Value = 5176600
$ws.Range("J9").Value = 4734500

# Row 10: Gross Profit
$ws.Range("D10").Value = 3228800
$ws.Range("E10").Value = 2582100
$ws.Range("F10").Value = 2425800
$ws.Range("G10").Value = 2186800
$ws.Range("H10").Value = 1809900
$ws.Range("I10").Value = 1235200
$ws.Range("J10").Value = 1433700

# Row 12: Research Development
$ws.Range("D12").Value = 501200
$ws.Range("E12").Value = 477400
$ws.Range("F12").Value = 469900
$ws.Range("G12").Value = 408400
$ws.Range("H12").Value = 341800
$ws.Range("I12").Value = 309900
$ws.Range("J12").Value = 271700

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 11936900
$ws.Range("E17").Value = 9581900
$ws.Range("F17").Value = 9588100
$ws.Range("G17").Value = 8293700
$ws.Range("H17").Value = 7143800
$ws.Range("I17").Value = 6252700
$ws.Range("J17").Value = 5507600

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = 1515400
$ws.Range("E18").Value = 1259900
$ws.Range("F18").Value = 1063700
$ws.Range("G18").Value = 1002900
$ws.Range("H18").Value = 767200
$ws.Range("I18").Value = 159100
$ws.Range("J18").Value = 660600

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = 29800
$ws.Range("E20").Value = 41900
$ws.Range("F20").Value = 16600
$ws.Range("G20").Value = -21300
$ws.Range("H20").Value = 10100
$ws.Range("I20").Value = -31800
$ws.Range("J20").Value = -17300

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 2164400
$ws.Range("E21").Value = 1842400
$ws.Range("F21").Value = 1668000
$ws.Range("G21").Value = 1389600
$ws.Range("H21").Value = 1134700
$ws.Range("I21").Value = 443400
$ws.Range("J21").Value = 928400

# Row 22: Interest Expense
$ws.Range("D22").Value = 58500
$ws.Range("E22").Value = 24300
$ws.Range("F22").Value = 21100
$ws.Range("G22").Value = 13400
$ws.Range("H22").Value = 13800
$ws.Range("I22").Value = 6100
$ws.Range("J22").Value = 2700

# Row 23: Income Before Tax
$ws.Range("D23").Value = 1486700
$ws.Range("E23").Value = 1277500
$ws.Range("F23").Value = 1059200
$ws.Range("G23").Value = 968100
$ws.Range("H23").Value = 763500
$ws.Range("I23").Value = 121100
$ws.Range("J23").Value = 640500

# Row 24: Income Tax Expense
$ws.Range("D24").Value = 292300
$ws.Range("E24").Value = 265400
$ws.Range("F24").Value = 236500
$ws.Range("G24").Value = 262500
$ws.Range("H24").Value = 231900
$ws.Range("I24").Value = 59300
$ws.Range("J24").Value = 170000

# Row 26: Income After Tax
$ws.Range("D26").Value = 1194400
$ws.Range("E26").Value = 1012100
$ws.Range("F26").Value = 822600
$ws.Range("G26").Value = 705700
$ws.Range("H26").Value = 531600
$ws.Range("I26").Value = 61800
$ws.Range("J26").Value = 470600

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = 1188200
$ws.Range("E27").Value = 1003500
$ws.Range("F27").Value = 813100
$ws.Range("G27").Value = 687200
$ws.Range("H27").Value = 508700
$ws.Range("I27").Value = 72200
$ws.Range("J27").Value = 438400

# Row 29: Discontinued Operations
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = -70200

# Row 32: Other Items
$ws.Range("D32").Value = -29800
$ws.Range("E32").Value = -41900
$ws.Range("F32").Value = -16600
$ws.Range("G32").Value = 21300
$ws.Range("H32").Value = -10100
$ws.Range("I32").Value = 31800
$ws.Range("J32").Value = 17300

# Row 33: Net Income
$ws.Range("D33").Value = 1188200
$ws.Range("E33").Value = 1003500
$ws.Range("F33").Value = 813100
$ws.Range("G33").Value = 687200
$ws.Range("H33").Value = 508700
$ws.Range("I33").Value = 72200
$ws.Range("J33").Value = 368200

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = 1188200
$ws.Range("E35").Value = 1003500
$ws.Range("F35").Value = 813100
$ws.Range("G35").Value = 687200
$ws.Range("H35").Value = 508700
$ws.Range("I35").Value = 72200
$ws.Range("J35").Value = 368200

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 4808300
$ws.Range("E41").Value = 2907100
$ws.Range("F41").Value = 5531400
$ws.Range("G41").Value = 2439900
$ws.Range("H41").Value = 2239600
$ws.Range("I41").Value = 1748500
$ws.Range("J41").Value = 1177800

# Row 42: Short Term Investments
$ws.Range("D42").Value = 15500
$ws.Range("E42").Value = 26700
$ws.Range("F42").Value = 18200

# Row 43: Net Receivables
$ws.Range("D43").Value = 7050200
$ws.Range("E43").Value = 3169200
$ws.Range("F43").Value = 2453200
$ws.Range("G43").Value = 2328000
$ws.Range("H43").Value = 1900100
$ws.Range("I43").Value = 1438100
$ws.Range("J43").Value = 1649500

# Row 44: Inventory
$ws.Range("D44").Value = 4224700
$ws.Range("E44").Value = 1905900
$ws.Range("F44").Value = 1661700
$ws.Range("G44").Value = 1987800
$ws.Range("H44").Value = 1471900
$ws.Range("I44").Value = 902400
$ws.Range("J44").Value = 826700

# Row 45: Other Current Assets
$ws.Range("D45").Value = 450600
$ws.Range("E45").Value = 129200
$ws.Range("F45").Value = 172500
$ws.Range("G45").Value = 190000
$ws.Range("H45").Value = 656800
$ws.Range("I45").Value = 437200
$ws.Range("J45").Value = 317100

# Row 46: Total Current Assets
$ws.Range("D46").Value = 8287200
$ws.Range("E46").Value = 8138100
$ws.Range("F46").Value = 6826000
$ws.Range("G46").Value = 6590300
$ws.Range("H46").Value = 5572400
$ws.Range("I46").Value = 4526200
$ws.Range("J46").Value = 3971200

# Row 47: Long Term Investments
$ws.Range("D47").Value = 462500
$ws.Range("E47").Value = 221200
$ws.Range("F47").Value = 339900
$ws.Range("G47").Value = 214100
$ws.Range("H47").Value = 166800
$ws.Range("I47").Value = 154200
$ws.Range("J47").Value = 140800

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 8123500
$ws.Range("E48").Value = 3613600
$ws.Range("F48").Value = 483100
$ws.Range("G48").Value = 3092800
$ws.Range("H48").Value = 3402300
$ws.Range("I48").Value = 3799100
$ws.Range("J48").Value = 2395800

# Row 49: Goodwill
$ws.Range("D49").Value = 4351500
$ws.Range("E49").Value = 3072200
$ws.Range("F49").Value = 3634600
$ws.Range("G49").Value = 2141700
$ws.Range("H49").Value = 2027100
$ws.Range("I49").Value = 1200300
$ws.Range("J49").Value = 727900

# Row 52: Other Assets
$ws.Range("D52").Value = 221600
$ws.Range("E52").Value = 133000
$ws.Range("F52").Value = 961800
$ws.Range("G52").Value = 1040800
$ws.Range("H52").Value = 82800

# Row 54: Total Assets
$ws.Range("D54").Value = 16037000
$ws.Range("E54").Value = 15178100
$ws.Range("F54").Value = 12444800
$ws.Range("G54").Value = 12270400
$ws.Range("H54").Value = 10549100
$ws.Range("I54").Value = 9089000
$ws.Range("J54").Value = 7235600

# Row 57: Accounts Payable
$ws.Range("D57").Value = 5716800
$ws.Range("E57").Value = 2262700
$ws.Range("F57").Value = 3289300
$ws.Range("G57").Value = 1762800
$ws.Range("H57").Value = 1504100
$ws.Range("I57").Value = 1212900
$ws.Range("J57").Value = 970400

# Row 58: Short/Current Long Term Debt
$ws.Range("D58").Value = 549000
$ws.Range("E58").Value = 2265800
$ws.Range("F58").Value = 1481400
$ws.Range("G58").Value = 884900
$ws.Range("H58").Value = 468700
$ws.Range("I58").Value = 1504500
$ws.Range("J58").Value = 789000

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 1511800
$ws.Range("E59").Value = 886700
$ws.Range("F59").Value = 1026700
$ws.Range("G59").Value = 633400
$ws.Range("H59").Value = 581500
$ws.Range("I59").Value = 581100
$ws.Range("J59").Value = 521900

# Row 60: Total Current Liabilities
$ws.Range("D60").Value = 4086200
$ws.Range("E60").Value = 5415200
$ws.Range("F60").Value = 3793700
$ws.Range("G60").Value = 3281000
$ws.Range("H60").Value = 2554300
$ws.Range("I60").Value = 3298500
$ws.Range("J60").Value = 2281300

# Row 61: Long Term Debt
$ws.Range("D61").Value = 2844300
$ws.Range("E61").Value = 1462500
$ws.Range("F61").Value = 1236700
$ws.Range("G61").Value = 1668900
$ws.Range("H61").Value = 2706700
$ws.Range("I61").Value = 1322300
$ws.Range("J61").Value = 915200

# Row 62: Other Liabilities
$ws.Range("D62").Value = 587400
$ws.Range("E62").Value = 563900
$ws.Range("F62").Value = 441200
$ws.Range("G62").Value = 512400
$ws.Range("H62").Value = 399400
$ws.Range("I62").Value = 365700
$ws.Range("J62").Value = 191600

# Row 66: Total Liabilities
$ws.Range("D66").Value = 7607200
$ws.Range("E66").Value = 7525100
$ws.Range("F66").Value = 5547100
$ws.Range("G66").Value = 5535800
$ws.Range("H66").Value = 5866700
$ws.Range("I66").Value = 5331500
$ws.Range("J66").Value = 3889200

# Row 72: Retained Earnings
$ws.Range("D72").Value = 7436200
$ws.Range("E72").Value = 6471800
$ws.Range("F72").Value = 5651500
$ws.Range("G72").Value = 3865900
$ws.Range("H72").Value = 3322100
$ws.Range("I72").Value = 2916600
$ws.Range("J72").Value = 2954100

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 8429800
$ws.Range("E76").Value = 7653000
$ws.Range("F76").Value = 6897700
$ws.Range("G76").Value = 6734500
$ws.Range("H76").Value = 4682500
$ws.Range("I76").Value = 3757500
$ws.Range("J76").Value = 3346400

# Row 81: Net Income
$ws.Range("D81").Value = 1188200
$ws.Range("E81").Value = 1003500
$ws.Range("F81").Value = 813100
$ws.Range("G81").Value = 687200
$ws.Range("H81").Value = 508700
$ws.Range("I81").Value = 72200
$ws.Range("J81").Value = 368200

# Row 83: Depreciation
$ws.Range("D83").Value = 618500
$ws.Range("E83").Value = 540000
$ws.Range("F83").Value = 587100
$ws.Range("G83").Value = 407600
$ws.Range("H83").Value = 356900
$ws.Range("I83").Value = 315800
$ws.Range("J83").Value = 284900

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 1587100
$ws.Range("E89").Value = 1173900
$ws.Range("F89").Value = 1334800
$ws.Range("G89").Value = 830600
$ws.Range("H89").Value = 788500
$ws.Range("I89").Value = 997000
$ws.Range("J89").Value = 512700

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -821200
$ws.Range("E91").Value = -621200
$ws.Range("F91").Value = -1480900
$ws.Range("G91").Value = -524700
$ws.Range("H91").Value = -364300
$ws.Range("I91").Value = -554800
$ws.Range("J91").Value = -374700

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -1029800
$ws.Range("E94").Value = -1911700
$ws.Range("F94").Value = -862200
$ws.Range("G94").Value = -734300
$ws.Range("H94").Value = -571100
$ws.Range("I94").Value = -1210000
$ws.Range("J94").Value = -180100

# Row 96: Dividends Paid
$ws.Range("D96").Value = -241100
$ws.Range("E96").Value = -214500
$ws.Range("F96").Value = -214200
$ws.Range("G96").Value = -143400
$ws.Range("H96").Value = -103300
$ws.Range("I96").Value = -109600
$ws.Range("J96").Value = -112100

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -1056400
$ws.Range("E100").Value = 866500
$ws.Range("F100").Value = 70300
$ws.Range("G100").Value = -176400
$ws.Range("H100").Value = 121800
$ws.Range("I100").Value = 552500
$ws.Range("J100").Value = -7400

# Row 101: Effect Of Exchange Rate Changes 
$ws.Range("D101").Value = -3900
$ws.Range("E101").Value = 12800
$ws.Range("F101").Value = -217100
$ws.Range("G101").Value = 280500
$ws.Range("H101").Value = 151900
$ws.Range("I101").Value = 231300
$ws.Range("J101").Value = -100

# Row 102: Change In Cash and Cash Equivalents 
$ws.Range("D102").Value = -502900
$ws.Range("E102").Value = 141400
$ws.Range("F102").Value = 325800
$ws.Range("G102").Value = 200300
$ws.Range("H102").Value = 491100
$ws.Range("I102").Value = 570700
$ws.Range("J102").Value = 325200
